$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5764
    3  = 6330
    4  = 3300
    5  = 18174
    6  = 15489
    7  = 1406
    8  = 4296
    9  = 8026
    10 = 3274
    11 = 18319
    12 = 13587
    13 = 9147
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
